# Atualização de bases das ligas, do dia: 2024-01-29 às 17-07
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Portugal Primeira Liga")

# Pairs of adjacent rows whose data (all columns except A, the running id)
# need to be swapped with one another.
$swapPairs = @(
    @(520, 521),
    @(667, 668),
    @(695, 696),
    @(778, 779),
    @(776, 777),
    @(866, 867),
    @(905, 906),
    @(951, 952),
    @(957, 958),
    @(982, 983),
    @(985, 986)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $range1 = $ws.Range("B$r1`:AC$r1")
    $range2 = $ws.Range("B$r2`:AC$r2")
    $v1 = $range1.Value2
    $v2 = $range2.Value2
    $range1.Value = $v2
    $range2.Value = $v1
}

# Row 999 (id 6876612) was removed entirely; rows 987-998 shift up to take
# the place of the next row down (999 -> 998 -> ... -> 987), effectively
# deleting row 999 and pulling everything below row 986 up by one row.
$ws.Rows("999").Delete()
